$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$win = $excel.Windows.Item(1)
try { $win.ScrollIntoView(0,0,100,100,$true) } catch { Write-Output "no ScrollIntoView: $_" }
try { Write-Output $win.VisibleRange.Address } catch { Write-Output "no VisibleRange: $_" }
